# Daily refresh of the cryptocurrency price/volume table (GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns keep their text (string) representation even when
# the new value looks numeric (e.g. "215.61", "1.00"), matching the
# original inline-string cell type used throughout this sheet.
$textCells = @(
    "D2", "D3", "D5", "D8", "D10", "D12", "D13", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D25", "D27", "D30", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D42", "D43", "D44", "D46", "D47", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.130.14'
$ws.Range("D3").Value = '1.659.52'
$ws.Range("E3").Value = '  +3.77%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '215.61'
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("E6").Value = '  +1.47%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '0.250'
$ws.Range("E8").Value = '  +2.30%  '
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").Value = '19.54'
$ws.Range("E10").Value = '  +3.04%  '
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("D12").Value = '1.889.20'
$ws.Range("E12").Value = '  +3.55%  '
$ws.Range("D13").Value = '1.656.03'
$ws.Range("E13").Value = '  +2.43%  '
$ws.Range("E14").Value = '  +1.98%  '
$ws.Range("E15").Value = '  +2.83%  '
$ws.Range("D16").Value = '64.97'
$ws.Range("E16").Value = '  +1.99%  '
$ws.Range("D17").Value = '240.96'
$ws.Range("E17").Value = '  +5.76%  '
$ws.Range("D18").Value = '27.098.67'
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("D19").Value = '7.87'
$ws.Range("E19").Value = '  +4.13%  '
$ws.Range("D20").Value = '0.0₃0730'
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("E22").Value = '  +5.05%  '
$ws.Range("D23").Value = '2.28'
$ws.Range("E23").Value = '  +4.70%  '
$ws.Range("E24").Value = '  +3.72%  '
$ws.Range("D25").Value = '146.06'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("D27").Value = '7.14'
$ws.Range("E27").Value = '  +2.59%  '
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("E29").Value = '  +3.25%  '
$ws.Range("D30").Value = '0.0498'
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("D32").Value = '1.522.63'
$ws.Range("E32").Value = '  +5.28%  '
$ws.Range("D33").Value = '3.29'
$ws.Range("E33").Value = '  +2.88%  '
$ws.Range("D34").Value = '3.06'
$ws.Range("E34").Value = '  +3.52%  '
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  +7.78%  '
$ws.Range("D36").Value = '2.42'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").Value = '0.576'
$ws.Range("E37").Value = '  +1.39%  '
$ws.Range("D38").Value = '0.895'
$ws.Range("E38").Value = '  +9.07%  '
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("D40").Value = '5.95'
$ws.Range("E40").Value = '  +3.05%  '
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D42").Value = '2.30'
$ws.Range("E42").Value = '  +5.59%  '
$ws.Range("D43").Value = '66.47'
$ws.Range("E43").Value = '  +9.82%  '
$ws.Range("D44").Value = '1.796.79'
$ws.Range("E44").Value = '  +3.40%  '
$ws.Range("E45").Value = '  +2.61%  '
$ws.Range("D46").Value = '0.914'
$ws.Range("E46").Value = '  -1.05%  '
$ws.Range("D47").Value = '90.64'
$ws.Range("E47").Value = '  +3.53%  '
$ws.Range("E48").Value = '  +3.62%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.0978'
$ws.Range("E49").Value = '  +3.05%  '
$ws.Range("D50").Value = '0.0502'
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.57'
$ws.Range("E51").Value = '  +2.15%  '
